# Update "想去人数" (attendance/interest count) figures in the "展览" sheet
# and the corresponding rows in the merged "全部类型" sheet, matching the
# freshly regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value2  = 250
$wsExhibit.Range("F3").Value2  = 2479
$wsExhibit.Range("F5").Value2  = 900
$wsExhibit.Range("F6").Value2  = 32
$wsExhibit.Range("F7").Value2  = 1306
$wsExhibit.Range("F8").Value2  = 1654
$wsExhibit.Range("F9").Value2  = 172
$wsExhibit.Range("F11").Value2 = 2318
$wsExhibit.Range("F12").Value2 = 474
$wsExhibit.Range("F13").Value2 = 146
$wsExhibit.Range("F14").Value2 = 52
$wsExhibit.Range("F16").Value2 = 102
$wsExhibit.Range("F17").Value2 = 91
$wsExhibit.Range("F18").Value2 = 8437
$wsExhibit.Range("F20").Value2 = 6498
$wsExhibit.Range("F21").Value2 = 10487
$wsExhibit.Range("F23").Value2 = 182
$wsExhibit.Range("F24").Value2 = 199
$wsExhibit.Range("F25").Value2 = 286
$wsExhibit.Range("F26").Value2 = 516
$wsExhibit.Range("F28").Value2 = 163
$wsExhibit.Range("F29").Value2 = 64
$wsExhibit.Range("F30").Value2 = 9
$wsExhibit.Range("F31").Value2 = 8
$wsExhibit.Range("F32").Value2 = 4430
$wsExhibit.Range("F33").Value2 = 310
$wsExhibit.Range("F34").Value2 = 419

# --- Sheet 4: 全部类型 (All types, merged view) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value2  = 250
$wsAll.Range("F6").Value2  = 2479
$wsAll.Range("F8").Value2  = 900
$wsAll.Range("F9").Value2  = 32
$wsAll.Range("F10").Value2 = 1306
$wsAll.Range("F12").Value2 = 1654
$wsAll.Range("F14").Value2 = 172
$wsAll.Range("F15").Value2 = 2318
$wsAll.Range("F17").Value2 = 474
$wsAll.Range("F18").Value2 = 146
$wsAll.Range("F19").Value2 = 52
$wsAll.Range("F22").Value2 = 102
$wsAll.Range("F23").Value2 = 91
$wsAll.Range("F24").Value2 = 8438
$wsAll.Range("F26").Value2 = 6498
$wsAll.Range("F27").Value2 = 10487
$wsAll.Range("F30").Value2 = 182
$wsAll.Range("F31").Value2 = 199
$wsAll.Range("F32").Value2 = 286
$wsAll.Range("F34").Value2 = 516
$wsAll.Range("F39").Value2 = 163
$wsAll.Range("F40").Value2 = 4430
$wsAll.Range("F47").Value2 = 419
